$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,4
$data[0,0] = [double]"0.071432925760746"
$data[0,1] = [double]"0.9800567626953125"
$data[0,2] = [double]"0.01750546507537365"
$data[0,3] = [double]"0.9988094568252563"
$data[1,0] = [double]"0.01207085885107517"
$data[1,1] = [double]"0.9981154203414917"
$data[1,2] = [double]"0.00889816228300333"
$data[1,3] = [double]"0.998869001865387"
$data[2,0] = [double]"0.007044259458780289"
$data[2,1] = [double]"0.9984675049781799"
$data[2,2] = [double]"0.005572132766246796"
$data[2,3] = [double]"0.9994047284126282"
$data[3,0] = [double]"0.002788325306028128"
$data[3,1] = [double]"0.9993994235992432"
$data[3,2] = [double]"0.003888229373842478"
$data[3,3] = [double]"0.9997618794441223"
$data[4,0] = [double]"0.002024667337536812"
$data[4,1] = [double]"0.9994615316390991"
$data[4,2] = [double]"0.0008195702685043216"
$data[4,3] = [double]"0.9998809695243835"
$data[5,0] = [double]"0.001594877801835537"
$data[5,1] = [double]"0.9995858073234558"
$data[5,2] = [double]"0.002492917934432626"
$data[5,3] = [double]"0.9998214244842529"
$data[6,0] = [double]"0.0009395657107234001"
$data[6,1] = [double]"0.9996893405914307"
$data[6,2] = [double]"0.003666960867121816"
$data[6,3] = [double]"0.9996428489685059"
$data[7,0] = [double]"0.001343700452707708"
$data[7,1] = [double]"0.9997721910476685"
$data[7,2] = [double]"0.002518505323678255"
$data[7,3] = [double]"0.9997618794441223"
$data[8,0] = [double]"0.0007529841968789697"
$data[8,1] = [double]"0.9998757243156433"
$data[8,2] = [double]"0.001128750038333237"
$data[8,3] = [double]"0.9998809695243835"
$data[9,0] = [double]"0.0008054365171119571"
$data[9,1] = [double]"0.9998136162757874"
$data[9,2] = [double]"0.001376897096633911"
$data[9,3] = [double]"0.9998214244842529"
$data[10,0] = [double]"0.001962084788829088"
$data[10,1] = [double]"0.9996686577796936"
$data[10,2] = [double]"0.003031467320397496"
$data[10,3] = [double]"0.9997023344039917"
$data[11,0] = [double]"0.0005426937714219093"
$data[11,1] = [double]"0.9998342990875244"
$data[11,2] = [double]"0.006332992110401392"
$data[11,3] = [double]"0.9992261528968811"
$data[12,0] = [double]"0.0003821065474767238"
$data[12,1] = [double]"0.9999171495437622"
$data[12,2] = [double]"0.001758369733579457"
$data[12,3] = [double]"0.9997023344039917"
$data[13,0] = [double]"0.0006284148548729718"
$data[13,1] = [double]"0.9998964667320251"
$data[13,2] = [double]"0.001030823099426925"
$data[13,3] = [double]"0.9998214244842529"
$data[14,0] = [double]"0.0007504862733185291"
$data[14,1] = [double]"0.9998136162757874"
$data[14,2] = [double]"0.001786935608834028"
$data[14,3] = [double]"0.9997023344039917"
$data[15,0] = [double]"0.0002922247222159058"
$data[15,1] = [double]"0.9998757243156433"
$data[15,2] = [double]"0.001108429976738989"
$data[15,3] = [double]"0.9998214244842529"
$data[16,0] = [double]"0.0006888119387440383"
$data[16,1] = [double]"0.9998342990875244"
$data[16,2] = [double]"0.001016682712361217"
$data[16,3] = [double]"0.9999404549598694"
$data[17,0] = [double]"0.000399757525883615"
$data[17,1] = [double]"0.9999171495437622"
$data[17,2] = [double]"0.0009421196300536394"
$data[17,3] = [double]"0.9999404549598694"
$data[18,0] = [double]"0.0004215103108435869"
$data[18,1] = [double]"0.999937891960144"
$data[18,2] = [double]"0.002302622888237238"
$data[18,3] = [double]"0.9996428489685059"
$data[19,0] = [double]"0.001163725857622921"
$data[19,1] = [double]"0.9997721910476685"
$data[19,2] = [double]"0.001165747991763055"
$data[19,3] = [double]"0.9999404549598694"
$data[20,0] = [double]"0.0001637792156543583"
$data[20,1] = [double]"0.9999793171882629"
$data[20,2] = [double]"0.00134793680626899"
$data[20,3] = [double]"0.9999404549598694"
$data[21,0] = [double]"0.000131393302581273"
$data[21,1] = [double]"0.9999585747718811"
$data[21,2] = [double]"0.002220623660832644"
$data[21,3] = [double]"0.9996428489685059"
$data[22,0] = [double]"0.0003488792572170496"
$data[22,1] = [double]"0.9999171495437622"
$data[22,2] = [double]"0.001558560528792441"
$data[22,3] = [double]"0.9999404549598694"
$data[23,0] = [double]"0.0009291154565289617"
$data[23,1] = [double]"0.9998342990875244"
$data[23,2] = [double]"0.001481929793953896"
$data[23,3] = [double]"0.9998809695243835"
$data[24,0] = [double]"0.0001762874308042228"
$data[24,1] = [double]"0.9999585747718811"
$data[24,2] = [double]"0.001393180689774454"
$data[24,3] = [double]"0.9999404549598694"
$data[25,0] = [double]"0.0003804616862908006"
$data[25,1] = [double]"0.9999793171882629"
$data[25,2] = [double]"0.001544533763080835"
$data[25,3] = [double]"0.9998809695243835"
$data[26,0] = [double]"0.0004437203169800341"
$data[26,1] = [double]"0.9998964667320251"
$data[26,2] = [double]"0.001726457616314292"
$data[26,3] = [double]"0.9999404549598694"
$data[27,0] = [double]"0.0007796509889885783"
$data[27,1] = [double]"0.9998964667320251"
$data[27,2] = [double]"0.001530378009192646"
$data[27,3] = [double]"0.9999404549598694"
$data[28,0] = [double]"0.0004197605885565281"
$data[28,1] = [double]"0.9999793171882629"
$data[28,2] = [double]"0.001683302572928369"
$data[28,3] = [double]"0.9999404549598694"
$data[29,0] = [double]"0.0004489562998060137"
$data[29,1] = [double]"0.999937891960144"
$data[29,2] = [double]"0.00163915054872632"
$data[29,3] = [double]"0.9999404549598694"
$data[30,0] = [double]"0.0002899367827922106"
$data[30,1] = [double]"0.9999171495437622"
$data[30,2] = [double]"0.001633108360692859"
$data[30,3] = [double]"0.9998809695243835"
$data[31,0] = [double]"9.961012256098911E-05"
$data[31,1] = [double]"0.9999793171882629"
$data[31,2] = [double]"0.001490167574957013"
$data[31,3] = [double]"0.9999404549598694"
$data[32,0] = [double]"0.0006253143073990941"
$data[32,1] = [double]"0.9998964667320251"
$data[32,2] = [double]"0.01001349650323391"
$data[32,3] = [double]"0.9951782822608948"
$data[33,0] = [double]"0.0002751315187197179"
$data[33,1] = [double]"0.999937891960144"
$data[33,2] = [double]"0.0009371892083436251"
$data[33,3] = [double]"0.9996428489685059"
$data[34,0] = [double]"8.158569835359231E-05"
$data[34,1] = [double]"0.9999793171882629"
$data[34,2] = [double]"0.000882901658769697"
$data[34,3] = [double]"0.9999404549598694"
$data[35,0] = [double]"2.096043135679793E-05"
$data[35,1] = [double]"1"
$data[35,2] = [double]"0.001120846485719085"
$data[35,3] = [double]"0.9999404549598694"
$data[36,0] = [double]"0.0007030910346657038"
$data[36,1] = [double]"0.9998757243156433"
$data[36,2] = [double]"0.001450333627872169"
$data[36,3] = [double]"0.9999404549598694"
$data[37,0] = [double]"5.391412560129538E-05"
$data[37,1] = [double]"1"
$data[37,2] = [double]"0.00163026072550565"
$data[37,3] = [double]"0.9997618794441223"
$data[38,0] = [double]"9.411505016032606E-05"
$data[38,1] = [double]"0.9999793171882629"
$data[38,2] = [double]"0.001584045006893575"
$data[38,3] = [double]"0.9999404549598694"
$data[39,0] = [double]"0.0002890351752284914"
$data[39,1] = [double]"0.9998964667320251"
$data[39,2] = [double]"0.001141366199590266"
$data[39,3] = [double]"0.9999404549598694"
$data[40,0] = [double]"0.0003681188682094216"
$data[40,1] = [double]"0.9999171495437622"
$data[40,2] = [double]"0.001533857779577374"
$data[40,3] = [double]"0.9999404549598694"
$data[41,0] = [double]"0.001020055147819221"
$data[41,1] = [double]"0.9998964667320251"
$data[41,2] = [double]"0.001804807921871543"
$data[41,3] = [double]"0.9998214244842529"
$data[42,0] = [double]"8.673007687320933E-05"
$data[42,1] = [double]"0.999937891960144"
$data[42,2] = [double]"0.001860575983300805"
$data[42,3] = [double]"0.9998214244842529"
$data[43,0] = [double]"3.469051080173813E-05"
$data[43,1] = [double]"0.9999793171882629"
$data[43,2] = [double]"0.001567209721542895"
$data[43,3] = [double]"0.9999404549598694"
$data[44,0] = [double]"0.0005204507033340633"
$data[44,1] = [double]"0.999937891960144"
$data[44,2] = [double]"0.001962702954187989"
$data[44,3] = [double]"0.9999404549598694"
$data[45,0] = [double]"0.000127913590404205"
$data[45,1] = [double]"0.9999585747718811"
$data[45,2] = [double]"0.001791693503037095"
$data[45,3] = [double]"0.9998809695243835"
$data[46,0] = [double]"5.919684917898849E-05"
$data[46,1] = [double]"0.9999793171882629"
$data[46,2] = [double]"0.001811030204407871"
$data[46,3] = [double]"0.9999404549598694"
$data[47,0] = [double]"6.218387716216967E-05"
$data[47,1] = [double]"0.9999793171882629"
$data[47,2] = [double]"0.001810107263736427"
$data[47,3] = [double]"0.9999404549598694"
$data[48,0] = [double]"0.0008113188669085503"
$data[48,1] = [double]"0.9999171495437622"
$data[48,2] = [double]"0.002404500730335712"
$data[48,3] = [double]"0.9999404549598694"
$data[49,0] = [double]"0.0002090093184961006"
$data[49,1] = [double]"0.9999585747718811"
$data[49,2] = [double]"0.002451836364343762"
$data[49,3] = [double]"0.9999404549598694"

$range = $ws.Range("A2:D51")
$range.Value2 = $data
